$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.693.68'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.897.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.45'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4918'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.67%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2937'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06741'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.896.70'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.24'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07252'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '90.84'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6759'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.036'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.61%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.692.40'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007983'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.65%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.141.06'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.809'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '192.07'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +34.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.086'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +4.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.371'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.52'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +12.26%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.408'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.302'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09076'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.001'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05220'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7418'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.108'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.762'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01833'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.678'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.123'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9294'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4394'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.10'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.733'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.77%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1351'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.66%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.534'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05860'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.690'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.78%  '

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.66'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.71%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3927'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.417'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.47%  '
